# BB: Update tasks list
#
# Applies the OOXML diff to the "SocialFund 1.0 Specs" document:
#  - wraps several words/phrases with <w:proofErr> spell/grammar markers,
#    which requires splitting existing <w:r> runs into multiple runs
#  - rewrites/reorders several bullets in the "Further features" list
#  - removes a block of now-obsolete bullets
#  - relocates the _GoBack bookmark from the last bullet to the
#    "button "Back"" bullet
#
# Strategy: use Range.InsertXML (which REPLACES the target range's
# contents with the supplied OOXML) on a per-paragraph basis, walking
# paragraphs from the bottom of the document upward so that edits to
# later paragraphs never invalidate the Paragraphs collection indices
# of paragraphs still to be processed. The big contiguous run of
# obsolete bullets is removed with a single Range.Delete() call.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Change group name" (last bullet, index 43 before any edits below):
#    drop the trailing _GoBack bookmark (it moves up to the
#    button "Back" bullet instead).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(43)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Change group name</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 2) Remove the long obsolete block of bullets (indices 31-42):
#      Add validation of fields on feedback controls.
#      User account: Remember me
#      Edit group or client name
#      Bug: User can't add other users to not own group
#      Bug: Fix paging for groups and users
#      Issue with additional library for grid
#      Remove GroupBox on CreateGroup page
#      Change admin rights
#      Additional information about user: Email, Phone number, Address
#      Send notification with proposition for all group members
#      button "Back"
#      Small Forum
#    (their content is either gone for good, or re-appears earlier in
#    the list via the edits below)
# ---------------------------------------------------------------------
$pFirst = $d.Paragraphs.Item(31)
$pLast = $d.Paragraphs.Item(42)
$rng = $d.Range($pFirst.Range.Start, $pLast.Range.End)
$rng.Delete()

# ---------------------------------------------------------------------
# 3) "Add "Remove coins" buttons" (index 30) ->
#    _GoBack bookmark + 'button "Back"'
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(30)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>button "Back"</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 4) "Mail format body" (index 29) -> "Change admin rights"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(29)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Change admin rights</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 5) "Add Register link near Log In link" (index 28) ->
#    "Edit group or client name"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(28)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Edit group or client name</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 6) "Logs" (index 27, carries <w:lastRenderedPageBreak/>) ->
#    "User account: " + "Remember me" (two runs)
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(27)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">User account: </w:t></w:r><w:r><w:t>Remember me</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 7) "CAPTCHA" (index 26) -> "Mail format body"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(26)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Mail format body</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 8) "Security features and vulnerabilities. " (index 25) -> "Logs"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(25)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Logs</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 9) "Notification to any users. " + "Send to email any actions
#    connection to user. " (index 24, two runs) ->
#    single run "Security features and vulnerabilities. "
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(24)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Security features and vulnerabilities. </w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 10) UI Design: "Groups page: ..." (index 18) -> wrap "Groups" with
#     <w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(18)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>Groups</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> page: User can create several groups and review table/list of groups where he is the member. Created group should be marked “Owner”, other groups – “Member”.</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 11) UI: "UI will be based on IndexServer web site." (index 15) ->
#     wrap "IndexServer" with spellStart/spellEnd
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(15)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>UI will be base</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t xml:space="preserve"> on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IndexServer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> web site.</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 12) Requirements bullet: 'someone paid in: should be next
#     information "money/value",  "person".' (index 7) -> split off a
#     gramStart/gramEnd wrapped ',  "' run
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(7)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Account manager loges-in to web site and add amount of money that </w:t></w:r><w:r><w:t>someone paid in: should be next information “money/value”</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,  “</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>person”.</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 13) Goal: "SocialFund is the way to ea" (index 5) -> wrap "SocialFund"
#     with spellStart/spellEnd
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(5)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>SocialFund</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is the way to ea</w:t></w:r><w:r><w:t>sy control money by all members of group.</w:t></w:r><w:r><w:t xml:space="preserve"> It can be any groups: parents of one school class, company department, friends and so on.</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 14) Environment: ", Github" (index 3) -> wrap "Github" with
#     spellStart/spellEnd
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(3)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Microsoft Visual Studio 2012, ASP.NET MVC 4.0, Entity Framework</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Microsoft SQL Server</w:t></w:r><w:r><w:t xml:space="preserve"> 2008</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 15) Title: "SocialFund 1.0 Specs" (index 1) -> wrap "SocialFund" with
#     spellStart/spellEnd
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(1)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>SocialFund</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> 1.0 Specs</w:t></w:r></w:p>
'@
$p.Range.InsertXML($xml)

Write-Output "done"
